$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.272653698921204
$ws.Range("B1").Value = 4.629242420196533
$ws.Range("C1").Value = 3.978886127471924
$ws.Range("D1").Value = 1.457034111022949
$ws.Range("E1").Value = 0.968996524810791
